$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "29.731.11"
$ws.Range("E2").Value = "  -3.46%  "
$ws.Range("D3").Value = "2.098.68"
$ws.Range("E3").Value = "  -2.58%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'345.05"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.5140"
$ws.Range("E7").Value = "  -2.90%  "
$ws.Range("D8").Value = "'0.4410"
$ws.Range("E8").Value = "  -3.85%  "
$ws.Range("D9").Value = "'52.74"
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("D10").Value = "'0.09239"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("D12").Value = "'24.91"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "2.100.60"
$ws.Range("E13").Value = "  -2.20%  "
$ws.Range("D14").Value = "'8.281"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").Value = "'6.758"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "'99.58"
$ws.Range("E16").Value = "  -2.82%  "
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "'20.86"
$ws.Range("E19").Value = "  +6.03%  "
$ws.Range("D20").Value = "'0.06630"
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("D23").Value = "29.789.56"
$ws.Range("E24").Value = "  -2.46%  "
$ws.Range("D25").Value = "'2.319"
$ws.Range("E25").Value = "  -3.27%  "
$ws.Range("D26").Value = "2.349.11"
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("D27").Value = "'21.91"
$ws.Range("E28").Value = "  -4.36%  "
$ws.Range("D29").Value = "'161.93"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").Value = "'132.94"
$ws.Range("E30").Value = "  -3.45%  "
$ws.Range("D31").Value = "'1.134"
$ws.Range("E31").Value = "  -7.89%  "
$ws.Range("E32").Value = "  -3.30%  "
$ws.Range("D33").Value = "'1.657"
$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("D34").Value = "'6.183"
$ws.Range("E34").Value = "  -3.81%  "
$ws.Range("D35").Value = "'3.938"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").Value = "'10.49"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'6.049"
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("D38").Value = "'0.02570"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("D39").Value = "'0.06737"
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("D41").Value = "'0.6866"
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("D42").Value = "'0.2231"
$ws.Range("E42").Value = "  -4.94%  "
$ws.Range("D43").Value = "'1.304"
$ws.Range("E43").Value = "  +1.96%  "
$ws.Range("D44").Value = "'0.6648"
$ws.Range("E44").Value = "  +2.05%  "
$ws.Range("D45").Value = "'14.32"
$ws.Range("E45").Value = "  -3.58%  "
$ws.Range("D46").Value = "'2.318"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").Value = "'3.615"
$ws.Range("E47").Value = "  -3.98%  "
$ws.Range("E48").Value = "  -5.47%  "
$ws.Range("D49").Value = "'1.222"
$ws.Range("E49").Value = "  -3.29%  "
$ws.Range("D50").Value = "'0.3366"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").Value = "'82.36"
$ws.Range("E51").Value = "  -1.61%  "
